$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text updates per diff. Values that look like plain numbers are
# force-written as text (leading apostrophe, like typing into Excel) and
# then reset to the Normal style so no stray number format sticks to the cell.
$ws.Range('D2').Value = '67.414.47'
$ws.Range('D3').Value = '3.223.02'
$ws.Range('E3').Value = '  -1.60%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''578.69'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.37%  '
$ws.Range('D6').Value = '''183.97'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.24%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '''0.606'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.76%  '
$ws.Range('D9').Value = '3.221.21'
$ws.Range('E9').Value = '  -1.66%  '
$ws.Range('E10').Value = '  -3.02%  '
$ws.Range('E11').Value = '  -2.36%  '
$ws.Range('E12').Value = '  -1.53%  '
$ws.Range('D13').Value = '3.777.88'
$ws.Range('E13').Value = '  -1.78%  '
$ws.Range('E14').Value = '  +0.05%  '
$ws.Range('D15').Value = '''27.69'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -3.31%  '
$ws.Range('D16').Value = '67.466.30'
$ws.Range('E16').Value = '  -0.85%  '
$ws.Range('E17').Value = '  -1.97%  '
$ws.Range('D18').Value = '3.201.30'
$ws.Range('E18').Value = '  -2.12%  '
$ws.Range('E19').Value = '  -1.95%  '
$ws.Range('D20').Value = '''13.43'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.70%  '
$ws.Range('D21').Value = '''396.19'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +3.79%  '
$ws.Range('D22').Value = '''7.54'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -2.34%  '
$ws.Range('D24').Value = '''71.20'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.35%  '
$ws.Range('D25').Value = '''0.515'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.10%  '
$ws.Range('E26').Value = '  -2.91%  '
$ws.Range('E27').Value = '  -2.96%  '
$ws.Range('D28').Value = '''9.51'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -3.17%  '
$ws.Range('D29').Value = '''1.00'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('E30').Value = '  -2.39%  '
$ws.Range('E31').Value = '  -5.39%  '
$ws.Range('D32').Value = '''22.57'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.32%  '
$ws.Range('E33').Value = '  -4.17%  '
$ws.Range('D34').Value = '''0.998'
$ws.Range('D34').Style = "Normal"
$ws.Range('E35').Value = '  -2.59%  '
$ws.Range('D36').Value = '''160.32'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.21%  '
$ws.Range('E37').Value = '  -4.83%  '
$ws.Range('D38').Value = '''1.88'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.76%  '
$ws.Range('D39').Value = '''26.40'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.85%  '
$ws.Range('D40').Value = '''0.802'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -4.53%  '
$ws.Range('E41').Value = '  -1.25%  '
$ws.Range('E42').Value = '  -4.51%  '
$ws.Range('E43').Value = '  -6.05%  '
$ws.Range('E44').Value = '  -1.75%  '
$ws.Range('D45').Value = '''40.62'
$ws.Range('D45').Style = "Normal"
$ws.Range('D46').Value = '2.592.78'
$ws.Range('E46').Value = '  -2.32%  '
$ws.Range('D47').Value = '''24.50'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -3.75%  '
$ws.Range('D48').Value = '''333.13'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.89%  '
$ws.Range('E49').Value = '  -2.72%  '
$ws.Range('E50').Value = '  -0.28%  '
$ws.Range('E51').Value = '  -1.73%  '
